$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture style swatches from existing well-formatted cells (A3 bold/top, B3 normal-wrap, C3 red-wrap) ---
# Placed far away (row 1000) so the swatch never collides with the real content area while we rebuild it.
$ws.Range("A3").Copy() | Out-Null
$ws.Cells.Item(1000, 1).PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Cells.Item(1000, 2).PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Cells.Item(1000, 3).PasteSpecial(-4122) | Out-Null

# --- Step 2: wipe all existing rows of the table (this also discards all old row heights) ---
$ws.Rows("1:23").Delete()

# Swatch row shifted up by 23 rows because of the delete above.
$swatchRow = 1000 - 23

# --- Step 3: rebuild every row of the new table, stamping the correct style from the swatch row, then setting content ---

# Row 1
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Value = "Ementa atual:"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

# Row 2
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Value = "LOM3250"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Value = "LOM3250"

# Row 3
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = "Nome:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Value = " Trabalho de Graduação II"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = " Trabalho de Graduação II"

# Row 4
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "Name:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = "Graduation Monograph II"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Value = "Graduation Monograph II"

# Row 5
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = "2"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Value = "2"

# Row 6
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Value = "4"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").Value = "4"

# Row 7
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = "Carga horária:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Value = "150 h"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Value = "150 h"

# Row 8
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = "Ativação:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = "01/01/2023"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "01/01/2023"

# Row 9
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "EF-9"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value = "EF-9"

# Row 10
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "Objetivos:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = "O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = "O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico."
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = "Objectives:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Value = "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"

# Row 14
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# Row 15
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = "Programa resumido:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = "Elaborar uma monografia de Trabalho de Graduação sob a orientação de docente e apresentá-la perante uma banca de examinadores."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = "Elaborar uma monografia de Trabalho de Graduação sob a orientação de docente e apresentá-la perante uma banca de examinadores."
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "Short syllabus:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners."
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "Programa:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = "O programa da disciplina será constituído pelas seguintes etapas: 1) Propor no início do período letivo um plano de trabalho a ser avaliado por uma comissão de professores. 2) Elaborar a monografia cujo tema seja pertencente ao conteúdo programático do curso de Engenharia Física, podendo ser um tópico de interesse técnico ou científico, estudo de caso ou uma proposta de projeto. 3) Definição e divulgação da data de apresentação após a entrega da monografia com antecedência de, no mínimo, 15 dias úteis. 4) Definição da banca de examinadores, sendo constituída pelo professor orientador e por no mínimo dois professores convidados. 5) Apresentação e avaliação do TG. 6) Divulgação da avaliação. Em caso de aprovação, deverá ser feita a entrega do exemplar final da monografia (cópia impressa e eletrônica) com o de acordo do professor orientador."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = "O programa da disciplina será constituído pelas seguintes etapas: 1) Propor no início do período letivo um plano de trabalho a ser avaliado por uma comissão de professores. 2) Elaborar a monografia cujo tema seja pertencente ao conteúdo programático do curso de Engenharia Física, podendo ser um tópico de interesse técnico ou científico, estudo de caso ou uma proposta de projeto. 3) Definição e divulgação da data de apresentação após a entrega da monografia com antecedência de, no mínimo, 15 dias úteis. 4) Definição da banca de examinadores, sendo constituída pelo professor orientador e por no mínimo dois professores convidados. 5) Apresentação e avaliação do TG. 6) Divulgação da avaliação. Em caso de aprovação, deverá ser feita a entrega do exemplar final da monografia (cópia impressa e eletrônica) com o de acordo do professor orientador."
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "Syllabus:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor."
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Avaliação:"

# Row 20
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = "Método:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value = "Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de um plano de trabalho."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de um plano de trabalho."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = "Critério:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Value = "A critério da banca de avaliação poderá ser estabelecido um prazo para revisão e/ou correção da monografia."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = "A critério da banca de avaliação poderá ser estabelecido um prazo para revisão e/ou correção da monografia."
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "Bibliografia:"
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value = "A ser definida no plano de trabalho."
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = "A ser definida no plano de trabalho."
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Cells.Item($swatchRow, 1).Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Value = "Requisitos:"

# Row 25
$ws.Cells.Item($swatchRow, 2).Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Value = "LOM3267 -  Trabalho de Graduação I  (Requisito)`n"
$ws.Cells.Item($swatchRow, 3).Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null
$ws.Range("C25").Value = "LOM3267 -  Trabalho de Graduação I  (Requisito)`n"
$ws.Rows.Item(25).RowHeight = 30

# --- Step 4: remove the temporary swatch row ---
$ws.Cells.Item($swatchRow, 1).EntireRow.Delete() | Out-Null

Write-Output $ws.UsedRange.Address()
